# PRO-standarder (HL7 CDA) workbook update
# Commit: "Modulus opdateret til Modulus Social. Vena har fået status
#          godkendt for CPD-DK og XDS Metadata"
#
# For this particular file the substantive, scriptable edit is the
# "last updated" date stamp: the sheet is renamed from
# "Opdateret d. 02-12-2025" to "Opdateret d. 05-12-2025". Excel keeps the
# workbook-scoped defined name (PRO_standarder__HL7_CDA_) in sync
# automatically because it references the sheet by name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Opdateret d. 05-12-2025"
